$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (160, 161) to the feed logs sheet
$ws.Cells.Item(160, 1).Value = 159
$ws.Cells.Item(160, 2).Value = 1
$ws.Cells.Item(160, 3).Value = "2024-06-18 08:17:12"
$ws.Cells.Item(160, 4).Value = 200
$ws.Cells.Item(160, 5).Value = 8

$ws.Cells.Item(161, 1).Value = 160
$ws.Cells.Item(161, 2).Value = 2
$ws.Cells.Item(161, 3).Value = "2024-06-18 08:17:12"
$ws.Cells.Item(161, 4).Value = 200
$ws.Cells.Item(161, 5).Value = 1
